$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# df_filters sheet: add the "aggsmall" filter example tables
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("df_filters")

# Pre-format the new block (rows 77-99, columns G:K) with the same
# "grey fill" style used by the existing example tables above it.
$styleSrc = $ws.Range("G70:K70")
$styleSrc.Copy($ws.Range("G77:K99"))

# --- Block 1: {{ df | aggsmall(1, 3) }} -----------------------------------
$ws.Cells.Item(79, 1).Value = "{{ df | aggsmall(1, 3) }}"
$ws.Cells.Item(79, 8).Value = "name"
$ws.Cells.Item(79, 9).Value = "b"
$ws.Cells.Item(79, 10).Value = "c"
$ws.Cells.Item(79, 11).Value = "d"

$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = "a"
$ws.Cells.Item(80, 9).Value = 4
$ws.Cells.Item(80, 10).Value = 1
$ws.Cells.Item(80, 11).Value = 1

$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = "b"
$ws.Cells.Item(81, 9).Value = 2
$ws.Cells.Item(81, 10).Value = 2
$ws.Cells.Item(81, 11).Value = 1

$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = "c"
$ws.Cells.Item(82, 9).Value = 6
$ws.Cells.Item(82, 10).Value = 5
$ws.Cells.Item(82, 11).Value = 1

$ws.Cells.Item(83, 7).Value = 3
$ws.Cells.Item(83, 8).Value = "d"
$ws.Cells.Item(83, 9).Value = 6
$ws.Cells.Item(83, 10).Value = 7
$ws.Cells.Item(83, 11).Value = 6

$ws.Cells.Item(84, 7).Value = 4
$ws.Cells.Item(84, 8).Value = "e"
$ws.Cells.Item(84, 9).Value = 9
$ws.Cells.Item(84, 10).Value = 8
$ws.Cells.Item(84, 11).Value = 7

# --- Block 2: {{ df }} ------------------------------------------------------
$ws.Cells.Item(86, 1).Value = "{{ df }}"
$ws.Cells.Item(86, 8).Value = "name"
$ws.Cells.Item(86, 9).Value = "b"
$ws.Cells.Item(86, 10).Value = "c"
$ws.Cells.Item(86, 11).Value = "d"

$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = "a"
$ws.Cells.Item(87, 9).Value = 4
$ws.Cells.Item(87, 10).Value = 1
$ws.Cells.Item(87, 11).Value = 1

$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = "b"
$ws.Cells.Item(88, 9).Value = 2
$ws.Cells.Item(88, 10).Value = 2
$ws.Cells.Item(88, 11).Value = 1

$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = "c"
$ws.Cells.Item(89, 9).Value = 6
$ws.Cells.Item(89, 10).Value = 5
$ws.Cells.Item(89, 11).Value = 1

$ws.Cells.Item(90, 7).Value = 3
$ws.Cells.Item(90, 8).Value = "d"
$ws.Cells.Item(90, 9).Value = 6
$ws.Cells.Item(90, 10).Value = 7
$ws.Cells.Item(90, 11).Value = 6

$ws.Cells.Item(91, 7).Value = 4
$ws.Cells.Item(91, 8).Value = "e"
$ws.Cells.Item(91, 9).Value = 9
$ws.Cells.Item(91, 10).Value = 8
$ws.Cells.Item(91, 11).Value = 7

# --- Block 3: {{ df | aggsmall(8, 3) }} -------------------------------------
$ws.Cells.Item(93, 1).Value = "{{ df | aggsmall(8, 3) }}"
$ws.Cells.Item(93, 8).Value = "name"
$ws.Cells.Item(93, 9).Value = "b"
$ws.Cells.Item(93, 10).Value = "c"
$ws.Cells.Item(93, 11).Value = "d"

$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = "Other"
$ws.Cells.Item(94, 9).Value = 27
$ws.Cells.Item(94, 10).Value = 23
$ws.Cells.Item(94, 11).Value = 16

# --- Block 4: {{ df | aggsmall(5, 3) }} -------------------------------------
$ws.Cells.Item(96, 1).Value = "{{ df | aggsmall(5, 3) }}"
$ws.Cells.Item(96, 8).Value = "name"
$ws.Cells.Item(96, 9).Value = "b"
$ws.Cells.Item(96, 10).Value = "c"
$ws.Cells.Item(96, 11).Value = "d"

$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = "d"
$ws.Cells.Item(97, 9).Value = 6
$ws.Cells.Item(97, 10).Value = 7
$ws.Cells.Item(97, 11).Value = 6

$ws.Cells.Item(98, 7).Value = 1
$ws.Cells.Item(98, 8).Value = "e"
$ws.Cells.Item(98, 9).Value = 9
$ws.Cells.Item(98, 10).Value = 8
$ws.Cells.Item(98, 11).Value = 7

$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = "Other"
$ws.Cells.Item(99, 9).Value = 12
$ws.Cells.Item(99, 10).Value = 8
$ws.Cells.Item(99, 11).Value = 3

# Update the sheet's remembered selection (this also momentarily makes
# df_filters the active sheet/tab).
$ws.Range("D73").Select()

# ------------------------------------------------------------------
# Restore the originally active tab: Sheet1 (instead of pic_filters)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
